$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns F (BASE AMOUNT), G (INITIAL AMOUNT), H (TOTAL)
# for rows 2 through 22, per the "code studies commit 24/11/2020" update.
$updates = @{
    2  = @(2000, 5000, 7000)
    3  = @(3000, 0,    4000)
    4  = @(3000, 5000, 9000)
    5  = @(3000, 5000, 9000)
    6  = @(2000, 5000, 7000)
    7  = @(2000, 5000, 7000)
    8  = @(3000, 0,    4000)
    9  = @(2000, 0,    2000)
    10 = @(3000, 0,    3800)
    11 = @(2000, 0,    2000)
    12 = @(3000, 0,    4000)
    13 = @(2000, 0,    2000)
    14 = @(3000, 0,    4000)
    15 = @(2000, 0,    2000)
    16 = @(3000, 0,    4000)
    17 = @(2000, 0,    2000)
    18 = @(3000, 0,    4000)
    19 = @(2000, 0,    2000)
    20 = @(3000, 0,    4000)
    21 = @(2000, 0,    2000)
    22 = @(3000, 0,    4000)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
}
